$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.932.30"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.555.21"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.488"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.91"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.246"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "1.775.61"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").Value = "1.554.15"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("D16").Value = "26.922.64"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "218.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0468"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.23%  "
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").Value = "1.433.56"
$ws.Range("E33").Value = "  +4.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.979"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.17%  "
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.519"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.809"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("E42").Value = "  -0.77%  "
$ws.Range("E44").Value = "  +2.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("E46").Value = "  +1.93%  "
$ws.Range("D47").Value = "1.689.73"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0525"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.44%  "
$ws.Range("E50").Value = "  +3.33%  "
